$wb = $excel.ActiveWorkbook

# --- Sheet TABLE_1: header row (row 4) new date columns EL (142) / EM (143) ---
$ws1 = $wb.Worksheets.Item("TABLE_1")
$ws1.Cells.Item(4, 142).Value = "'09/01/2023"
$ws1.Cells.Item(4, 143).Value = "'10/01/2023"

# --- Sheet TABLE_2: header row (row 4) new date columns DZ (130) / EA (131) ---
$ws2 = $wb.Worksheets.Item("TABLE_2")
$ws2.Cells.Item(4, 130).Value = "'09/01/2023"
$ws2.Cells.Item(4, 131).Value = "'10/01/2023"

# --- TABLE_1 data rows 5-56: update EJ (only row5), EK, and add EL/EM ---
$t1 = @(
    @{ Row=5; EJ=9039; EK=9484.3; EL=10584.9; EM=10929 },
    @{ Row=6; EJ=$null; EK=165; EL=171.5; EM=174.1 },
    @{ Row=7; EJ=$null; EK=21.4; EL=27.2; EM=28.2 },
    @{ Row=8; EJ=$null; EK=185.5; EL=200.6; EM=202.7 },
    @{ Row=9; EJ=$null; EK=87; EL=100.6; EM=102.6 },
    @{ Row=10; EJ=$null; EK=1111; EL=1184.1; EM=1234.2 },
    @{ Row=11; EJ=$null; EK=206.3; EL=224.2; EM=231.9 },
    @{ Row=12; EJ=$null; EK=97.8; EL=123.5; EM=127.4 },
    @{ Row=13; EJ=$null; EK=33.2; EL=37.9; EM=39.5 },
    @{ Row=14; EJ=$null; EK=$null; EL=$null; EM=$null },
    @{ Row=15; EJ=$null; EK=459.2; EL=484; EM=491.2 },
    @{ Row=16; EJ=$null; EK=344.2; EL=350.4; EM=352.9 },
    @{ Row=17; EJ=$null; EK=$null; EL=$null; EM=$null },
    @{ Row=18; EJ=$null; EK=57.1; EL=64; EM=67.7 },
    @{ Row=19; EJ=$null; EK=384.8; EL=445; EM=456.5 },
    @{ Row=20; EJ=$null; EK=206.5; EL=229.5; EM=235.3 },
    @{ Row=21; EJ=$null; EK=119.2; EL=142.1; EM=147.9 },
    @{ Row=22; EJ=$null; EK=112.2; EL=132.9; EM=138.6 },
    @{ Row=23; EJ=$null; EK=150.9; EL=161.2; EM=163.4 },
    @{ Row=24; EJ=$null; EK=140.5; EL=150.1; EM=151.8 },
    @{ Row=25; EJ=$null; EK=39.5; EL=50; EM=51.5 },
    @{ Row=26; EJ=$null; EK=197.7; EL=219.8; EM=224.9 },
    @{ Row=27; EJ=$null; EK=207.5; EL=240; EM=247.9 },
    @{ Row=28; EJ=$null; EK=272.7; EL=324; EM=334.1 },
    @{ Row=29; EJ=$null; EK=164.3; EL=194; EM=207.3 },
    @{ Row=30; EJ=$null; EK=99.6; EL=104.6; EM=105 },
    @{ Row=31; EJ=$null; EK=$null; EL=$null; EM=$null },
    @{ Row=32; EJ=$null; EK=33.4; EL=40.7; EM=42.4 },
    @{ Row=33; EJ=$null; EK=81; EL=92.7; EM=96 },
    @{ Row=34; EJ=$null; EK=69; EL=82.7; EM=85.2 },
    @{ Row=35; EJ=$null; EK=33.1; EL=45.3; EM=46.6 },
    @{ Row=36; EJ=$null; EK=259.6; EL=303.8; EM=323.7 },
    @{ Row=37; EJ=$null; EK=71.6; EL=78.6; EM=81.4 },
    @{ Row=38; EJ=$null; EK=532.9; EL=651.1; EM=675.4 },
    @{ Row=39; EJ=$null; EK=286.3; EL=319; EM=327.3 },
    @{ Row=40; EJ=$null; EK=28.8; EL=37.9; EM=39.8 },
    @{ Row=41; EJ=$null; EK=349.8; EL=383.8; EM=394.5 },
    @{ Row=42; EJ=$null; EK=129.1; EL=149.6; EM=153.5 },
    @{ Row=43; EJ=$null; EK=112.5; EL=126.9; EM=137.6 },
    @{ Row=44; EJ=$null; EK=278; EL=319.4; EM=327.8 },
    @{ Row=45; EJ=$null; EK=23.2; EL=28.3; EM=29.7 },
    @{ Row=46; EJ=$null; EK=150.8; EL=161; EM=164.2 },
    @{ Row=47; EJ=$null; EK=29.9; EL=37.3; EM=38.7 },
    @{ Row=48; EJ=$null; EK=182.1; EL=206.7; EM=212.6 },
    @{ Row=49; EJ=$null; EK=1087.7; EL=1156.5; EM=1191.9 },
    @{ Row=50; EJ=$null; EK=113.5; EL=128.2; EM=132 },
    @{ Row=51; EJ=$null; EK=26.5; EL=32.8; EM=34.1 },
    @{ Row=52; EJ=$null; EK=292.2; EL=324.4; EM=328.6 },
    @{ Row=53; EJ=$null; EK=218.8; EL=228.8; EM=248.2 },
    @{ Row=54; EJ=$null; EK=54.9; EL=59.4; EM=61.3 },
    @{ Row=55; EJ=$null; EK=166.8; EL=207; EM=216.9 },
    @{ Row=56; EJ=$null; EK=24.5; EL=29.2; EM=30.3 }
)

foreach ($r in $t1) {
    if ($r.EJ -ne $null) { $ws1.Cells.Item($r.Row, 140).Value = $r.EJ }
    if ($r.EK -ne $null) { $ws1.Cells.Item($r.Row, 141).Value = $r.EK }
    if ($r.EL -ne $null) { $ws1.Cells.Item($r.Row, 142).Value = $r.EL }
    if ($r.EM -ne $null) { $ws1.Cells.Item($r.Row, 143).Value = $r.EM }
}

# --- TABLE_2 data rows 5-56: update DY, and add DZ/EA ---
$t2 = @(
    @{ Row=5; DY=2.36586761071115; DZ=2.35263402181482; EA=2.60719348811883 },
    @{ Row=6; DY=0.671140939597312; DZ=0.823045267489715; EA=0.694042799305951 },
    @{ Row=7; DY=0.0000000000000166014657887874; DZ=0.740740740740751; EA=0.714285714285724 },
    @{ Row=8; DY=3.92156862745098; DZ=4.75195822454308; EA=4.26954732510289 },
    @{ Row=9; DY=0.578034682080925; DZ=0.599999999999994; EA=0.293255131964821 },
    @{ Row=10; DY=2.3963133640553; DZ=2.05119365681289; EA=1.76451187335093 },
    @{ Row=11; DY=4.03429147755925; DZ=2.70270270270269; EA=3.75838926174497 },
    @{ Row=12; DY=-0.407331975560073; DZ=1.72981878088962; EA=2.08333333333334 },
    @{ Row=13; DY=1.52905198776758; DZ=3.83561643835616; EA=5.61497326203209 },
    @{ Row=14; DY=$null; DZ=$null; EA=$null },
    @{ Row=15; DY=1.39103554868625; DZ=1.55266470835081; EA=1.80310880829015 },
    @{ Row=16; DY=2.74626865671641; DZ=3.24101355332941; EA=3.09669880222026 },
    @{ Row=17; DY=$null; DZ=$null; EA=$null },
    @{ Row=18; DY=8.34914611005691; DZ=5.26315789473685; EA=6.61417322834646 },
    @{ Row=19; DY=3.13588850174216; DZ=3.65711623573259; EA=3.65576748410536 },
    @{ Row=20; DY=-1.10153256704981; DZ=2.40963855421686; EA=4.1611332447986 },
    @{ Row=21; DY=0.336700336700341; DZ=0.995024875621895; EA=1.37080191912269 },
    @{ Row=22; DY=1.63043478260869; DZ=1.45038167938932; EA=2.28782287822878 },
    @{ Row=23; DY=3.07377049180326; DZ=2.21940393151554; EA=2.31684408265499 },
    @{ Row=24; DY=2.25618631732168; DZ=3.1615120274914; EA=2.22222222222223 },
    @{ Row=25; DY=0.50890585241731; DZ=0.806451612903223; EA=1.17878192534381 },
    @{ Row=26; DY=4.21718502899313; DZ=3.28947368421053; EA=3.68833563854311 },
    @{ Row=27; DY=1.91552062868368; DZ=0.671140939597325; EA=0.731409995936604 },
    @{ Row=28; DY=5.1677593521018; DZ=5.43442889684347; EA=4.47154471544716 },
    @{ Row=29; DY=0.244051250762664; DZ=1.51753008895866; EA=1.76730486008838 },
    @{ Row=30; DY=2.25872689938399; DZ=0.965250965250951; EA=0.767754318618053 },
    @{ Row=31; DY=$null; DZ=$null; EA=$null },
    @{ Row=32; DY=-3.18840579710145; DZ=-4.46009389671361; EA=-2.97482837528605 },
    @{ Row=33; DY=1.3767209011264; DZ=1.75631174533479; EA=1.15911485774499 },
    @{ Row=34; DY=-0.862068965517233; DZ=4.68354430379747; EA=5.97014925373134 },
    @{ Row=35; DY=0.607902735562319; DZ=2.02702702702702; EA=1.96936542669582 },
    @{ Row=36; DY=4.08981555733762; DZ=3.89876880984951; EA=5.19987000324994 },
    @{ Row=37; DY=3.61794500723589; DZ=5.22088353413655; EA=7.67195767195769 },
    @{ Row=38; DY=0.661125802795618; DZ=3.25087218522042; EA=0.640739085084183 },
    @{ Row=39; DY=0.632688927943745; DZ=0.726239343227032; EA=1.48837209302326 },
    @{ Row=40; DY=0; DZ=3.55191256830602; EA=3.64583333333333 },
    @{ Row=41; DY=0.865051903114187; DZ=1.42706131078225; EA=1.15384615384615 },
    @{ Row=42; DY=2.46031746031746; DZ=3.74479889042996; EA=4.70668485675307 },
    @{ Row=43; DY=1.6260162601626; DZ=1.9277108433735; EA=1.54981549815498 },
    @{ Row=44; DY=2.0932794711715; DZ=-0.77663870767319; EA=-0.485731633272641 },
    @{ Row=45; DY=-2.1097046413502; DZ=1.4336917562724; EA=3.12500000000001 },
    @{ Row=46; DY=1.00468854655057; DZ=-0.124069478908182; EA=0.121951219512188 },
    @{ Row=47; DY=2.39726027397259; DZ=2.1917808219178; EA=1.8421052631579 },
    @{ Row=48; DY=0.552181115405853; DZ=0.632911392405055; EA=1.9664268585132 },
    @{ Row=49; DY=2.22744360902256; DZ=2.23656294200848; EA=1.86308862490387 },
    @{ Row=50; DY=5.28756957328386; DZ=4.56769983686785; EA=4.5958795562599 },
    @{ Row=51; DY=5.57768924302788; DZ=6.84039087947881; EA=6.23052959501558 },
    @{ Row=52; DY=3.54358610914245; DZ=3.74160537256156; EA=3.92156862745099 },
    @{ Row=53; DY=-1.30807397383851; DZ=-2.63829787234042; EA=-0.241157556270094 },
    @{ Row=54; DY=1.66666666666668; DZ=-0.834724540901491; EA=-0.648298217179912 },
    @{ Row=55; DY=-0.239234449760752; DZ=2.12136161815491; EA=2.26308345120227 },
    @{ Row=56; DY=-0.406504065040656; DZ=0.689655172413803; EA=0.664451827242522 }
)

foreach ($r in $t2) {
    if ($r.DY -ne $null) { $ws2.Cells.Item($r.Row, 129).Value = $r.DY }
    if ($r.DZ -ne $null) { $ws2.Cells.Item($r.Row, 130).Value = $r.DZ }
    if ($r.EA -ne $null) { $ws2.Cells.Item($r.Row, 131).Value = $r.EA }
}

